$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits (alpha_distance_range / beta_distance_range rows) ---
# Row 2 (alpha_distance_range): Min 5.6 -> 5.5, Max 10.4 -> 10.5
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 10.5

# Row 3 (beta_distance_range): Min 5.7 -> 5.5, Max 9.3000000000000007 -> 9.5
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.5

# Row 4 (ratio_threshold_range): Min 0.8 -> 0.7 (Max 1.4 unchanged)
$ws.Range("B4").Value = 0.7

# --- Column C got a little wider once the new numbers were entered
# (best-fit width grew from 5.5 to 5.875 "characters" wide) ---
$ws.Columns.Item(3).ColumnWidth = 5.14
